# Update "想去人数" (F column, people-interested count) and one
# "最低票价" (G5) value across the 展览 / 演出 / 全部类型 sheets, matching
# the refreshed scrape output committed at 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 272
$ws.Range("F4").Value  = 1812
$ws.Range("G5").Value  = 70
$ws.Range("F6").Value  = 556
$ws.Range("F7").Value  = 556
$ws.Range("F8").Value  = 5142
$ws.Range("F12").Value = 998
$ws.Range("F13").Value = 353
$ws.Range("F17").Value = 3036
$ws.Range("F19").Value = 114
$ws.Range("F22").Value = 104
$ws.Range("F24").Value = 956
$ws.Range("F25").Value = 330
$ws.Range("F27").Value = 3397
$ws.Range("F29").Value = 2664
$ws.Range("F30").Value = 271
$ws.Range("F31").Value = 1676
$ws.Range("F32").Value = 3860
$ws.Range("F34").Value = 908
$ws.Range("F35").Value = 445
$ws.Range("F36").Value = 1215
$ws.Range("F37").Value = 20
$ws.Range("F38").Value = 965
$ws.Range("F39").Value = 1225
$ws.Range("F40").Value = 43
$ws.Range("F41").Value = 961
$ws.Range("F42").Value = 625
$ws.Range("F43").Value = 440
$ws.Range("F44").Value = 384
$ws.Range("F45").Value = 302
$ws.Range("F46").Value = 3526

# ---- Sheet "演出" (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F23").Value = 28

# ---- Sheet "全部类型" (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 272
$ws.Range("F4").Value  = 1812
$ws.Range("G5").Value  = 70
$ws.Range("F6").Value  = 556
$ws.Range("F7").Value  = 556
$ws.Range("F8").Value  = 5142
$ws.Range("F12").Value = 353
$ws.Range("F14").Value = 3036
$ws.Range("F17").Value = 114
$ws.Range("F24").Value = 104
$ws.Range("F25").Value = 956
$ws.Range("F26").Value = 330
$ws.Range("F27").Value = 3397
$ws.Range("F31").Value = 2664
$ws.Range("F32").Value = 1676
$ws.Range("F33").Value = 3860
$ws.Range("F36").Value = 908
$ws.Range("F37").Value = 1215
$ws.Range("F38").Value = 20
$ws.Range("F39").Value = 965
$ws.Range("F41").Value = 1225
$ws.Range("F42").Value = 43
$ws.Range("F43").Value = 961
$ws.Range("F44").Value = 625
$ws.Range("F45").Value = 384
$ws.Range("F46").Value = 28
$ws.Range("F48").Value = 302
$ws.Range("F49").Value = 3526
